$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158. Excel shifts rows 158:189 down to 159:190,
# and the new blank row 158 inherits formatting from the row above it.
$ws.Rows(158).Insert()

# The freshly inserted row 158 is blank; populate it with the same record
# that is now sitting in row 159 (i.e. what used to be row 158 before the
# insert), then overwrite just the date (column D) with the new reading.
$ws.Range("A159:R159").Copy()
$ws.Range("A158").PasteSpecial()

$ws.Range("D158").Value = 44522
